$d = $word.ActiveDocument

# --- 1. Merge "In conclusion..." paragraph runs ---
# " work, as" + " " + "I" + " said " -> " work, as I said "
[void]$d.Content.Find.Execute(" work, as I said ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " work, as I said ", 2)

# " " + "I" + " was just passing...didn't like" -> " I was just passing...didn't like"
# (searched/replaced without retyping the apostrophe so AutoCorrect can't swap it
#  for a curly quote, and without crossing the gramEnd proofErr boundary)
[void]$d.Content.Find.Execute(" I was just passing", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " I was just passing", 2)

# --- 2. Replace the tail of the document: split "So" out of the final
#        paragraph (wrapped in gramStart/gramEnd), then append the two new
#        diary-entry paragraphs, keeping the _GoBack bookmark on the new
#        final paragraph. ---
$last = $d.Paragraphs.Last
$tailRange = $last.Range

$newXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>So</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> my issue again was just not understanding the API.</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:lastRenderedPageBreak/><w:t>07/03/2018</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:t>Today with the indispensable help of Swen the bot now builds at the choke correctly, though the plan needs editing and bit of function need fleshing out, the bot now does more of what its supposed to do.  I was correct in it being a build site/position issue,.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$tailRange.InsertXML($newXml)

# InsertXML leaves a stray empty paragraph behind (the old final mark of the
# story) after inserting the new block paragraphs - remove it by deleting
# its paragraph mark, which merges it away without touching visible text.
$trailing = $d.Paragraphs.Last
$start = $trailing.Range.Start
$end = $trailing.Range.End
$cleanupRange = $d.Range($start - 1, $end)
[void]$cleanupRange.Delete(1, 1)

Write-Host "Done"
